$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set every "No" entry in the Run Mode column (C) to "Yes".
for ($r = 2; $r -le 33; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq "No") {
        $cell.Value = "Yes"
    }
}

# Update the sheet view: zoom level and selection.
$ws.Application.ActiveWindow.Zoom = 71
$ws.Range("C2:C33").Select()
